$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.572.39"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.449.13"
$ws.Range("E3").Value = "  -0.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.56"
$ws.Range("E5").Value = "  -1.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.65"
$ws.Range("E6").Value = "  -0.60%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -1.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.110"
$ws.Range("E9").Value = "  -1.56%  "

# Row 10
$ws.Range("E10").Value = "  -0.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.13"
$ws.Range("E11").Value = "  -2.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("E12").Value = "  -1.65%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.47"
$ws.Range("E13").Value = "  -1.96%  "

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000172"
$ws.Range("E14").Value = "  -3.26%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.893.44"
$ws.Range("E15").Value = "  -0.94%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.409.52"
$ws.Range("E16").Value = "  -0.94%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.448.90"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.69"
$ws.Range("E18").Value = "  -5.42%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("E19").Value = "  -3.21%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.46"
$ws.Range("E20").Value = "  -2.68%  "

# Row 21
$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.21"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -0.44%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.83"
$ws.Range("E24").Value = "  +3.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.68"
$ws.Range("E25").Value = "  -2.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "642.97"
$ws.Range("E26").Value = "  -3.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.566.33"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0944"
$ws.Range("E29").Value = "  -4.46%  "

# Row 30
$ws.Range("E30").Value = "  -3.93%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.76"
$ws.Range("E31").Value = "  -3.53%  "

# Row 32
$ws.Range("E32").Value = "  -3.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.131"
$ws.Range("E33").Value = "  -2.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").Value = "  -3.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "151.84"
$ws.Range("E36").Value = "  -0.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.59"
$ws.Range("E37").Value = "  -4.03%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.363"
$ws.Range("E38").Value = "  -2.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.43"
$ws.Range("E39").Value = "  -1.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.26"
$ws.Range("E40").Value = "  -3.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  -3.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.69"
$ws.Range("E42").Value = "  -3.50%  "

# Row 43
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("E44").Value = "  -0.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "152.50"
$ws.Range("E45").Value = "  +1.32%  "

# Row 46
$ws.Range("E46").Value = "  +1.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.51"
$ws.Range("E47").Value = "  -2.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.600"
$ws.Range("E48").Value = "  -0.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.81"
$ws.Range("E49").Value = "  -4.33%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0500"
$ws.Range("E50").Value = "  -2.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0900"
$ws.Range("E51").Value = "  -2.09%  "
